$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 645, shifting existing rows 645:666 down to 646:667
$ws.Rows.Item(645).Insert()

# Populate the newly inserted row 645 with the new weekly record
$ws.Cells.Item(645, 1).Value = 3
$ws.Cells.Item(645, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(645, 3).Value = "Coquimbo"
$ws.Cells.Item(645, 4).Value = 45239
$ws.Cells.Item(645, 5).Value = 5
$ws.Cells.Item(645, 6).Value = 100112027
$ws.Cells.Item(645, 7).Value = "Melón"
$ws.Cells.Item(645, 8).Value = "Tuna"
$ws.Cells.Item(645, 9).Value = "Primera"
$ws.Cells.Item(645, 10).Value = 60
$ws.Cells.Item(645, 11).Value = 21000
$ws.Cells.Item(645, 12).Value = 21000
$ws.Cells.Item(645, 13).Value = 21000
$ws.Cells.Item(645, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(645, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(645, 16).Value = 1750
$ws.Cells.Item(645, 17).Value = 12
$ws.Cells.Item(645, 18).Value = "Hortaliza"
